$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.092.68'
$ws.Range("E2").Value = '  +3.25%  '
$ws.Range("D3").Value = '3.109.08'
$ws.Range("E3").Value = '  +1.25%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '524.28'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.76%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.17'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.52%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.41%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.439'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.39%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.40'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.49%  '
$ws.Range("E10").Value = '  +1.18%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.384'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.27%  '
$ws.Range("D12").Value = '3.632.72'
$ws.Range("E12").Value = '  +1.17%  '
$ws.Range("E13").Value = '  +1.50%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.26'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.43%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000167'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.03%  '
$ws.Range("D16").Value = '59.017.16'
$ws.Range("E16").Value = '  +3.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.23'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +5.27%  '
$ws.Range("D18").Value = '3.082.29'
$ws.Range("E18").Value = '  +0.47%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.07'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.42%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.26'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '339.78'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.99%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.510'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.70%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.93'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.172'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.50%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").Value = '0.0₃0928'
$ws.Range("E27").Value = '  -2.46%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.64'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.89%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.34'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.02%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.84'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.87%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.23'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.16'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '155.31'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.32%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.65'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.51%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.12'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '27.44'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.42%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.30'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.86%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0686'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.19%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.97'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.51%  '
$ws.Range("D40").Value = '3.140.87'
$ws.Range("E40").Value = '  +1.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.89'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.17%  '
$ws.Range("E42").Value = '  -0.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.663'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.89%  '
$ws.Range("D44").Value = '2.291.02'
$ws.Range("E44").Value = '  +2.42%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.44'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.76%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0258'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.35%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '21.04'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.69%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.959'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.46%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.02'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.749'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +9.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '261.04'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +10.87%  '
